# fix up fields update and order for running effeciency #98
$wb = $excel.ActiveWorkbook

# --- Update view state (selection) on existing sheets ---
$ws1 = $wb.Worksheets.Item("gc_fields_display")
$ws1.Activate()
$ws1.Range("A86").Select()

$ws2 = $wb.Worksheets.Item("gc_fields_uom")
$ws2.Activate()
$ws2.Range("B68").Select()

# --- Add the new gc_fields_order sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "gc_fields_order"

# Header row
$ws3.Range("A1").Value = "field"
$ws3.Range("B1").Value = "field_order"
$ws3.Range("C1").Value = "category"
$ws3.Range("D1").Value = "activity_type"

# Data rows
$ws3.Range("A2").Value = "__CalcRunningEffectiveness"
$ws3.Range("B2").Value = 5
$ws3.Range("C2").Value = "run"

$ws3.Range("A3").Value = "__CalcMetabolicEfficiency"
$ws3.Range("B3").Value = 10
$ws3.Range("C3").Value = "run"
$ws3.Range("D3").Value = "running"

# Make new sheet active, matching the target selection/view state
$ws3.Activate()
$ws3.Range("A4").Select()
